# "Badge number search, other bug fixed"
#
# - Column A's header changes from a "sr #" counter to a "badge_nubmer"
#   (sic) column used for badge-number search.
# - The other header labels are normalized to snake_case / shortened:
#     name          (unchanged text, now in column B instead of D)
#     mobile number -> mobile_number
#     email address -> email
# - Column A gets a sensible width now that it holds real data.
# - A stray formatted (but empty) cell D2 is left behind with the
#   built-in "Hyperlink" style (looks like a hyperlink was added then
#   its contents cleared) and the sheet grows to A1:D2.
# - The active selection moves to D1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row text ---------------------------------------------------
$ws.Range("A1").Value = "badge_nubmer"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "mobile_number"
$ws.Range("D1").Value = "email"

# --- Column widths -------------------------------------------------------
# (stored OOXML "width" = ColumnWidth + 0.8333.., i.e. Excel's usual
# character-width padding; the numbers below land on the nearest
# reproducible width to the authored file.)
$ws.Columns.Item(1).ColumnWidth = 16.166666666666668   # -> stored width 17
$ws.Columns.Item(2).ColumnWidth = 33.33333333333333    # -> stored width ~34.17
$ws.Columns.Item(3).ColumnWidth = 14.0                 # -> stored width ~14.83
$ws.Columns.Item(4).ColumnWidth = 28.666666666666668   # -> stored width 29.5

# --- Leftover empty "Hyperlink"-styled cell at D2 -------------------------
# Adding then deleting a real hyperlink is what actually produces Excel's
# built-in Hyperlink cell style/font combo (with a proper theme color
# reference), leaving D2 formatted but empty and growing the sheet to
# A1:D2 - matching what a user gets after typing then clearing a link.
$ws.Hyperlinks.Add($ws.Range("D2"), "", "", "", "") | Out-Null
$ws.Hyperlinks.Delete()

# --- Selection -----------------------------------------------------------
$ws.Range("D1").Select()
